# carbon_stock_cfs/class_descriptions_key.xlsx
# "fixed typo in carbon stock classification key file, further edits to
#  all three tabs..."
#
# Content change: the NLCD class description "Glassland/Herbaceous" (row 14,
# column B) was a typo for "Grassland/Herbaceous". Fixing the text here is
# what actually reshuffles the shared-string table on save (Excel drops the
# now-unused string and appends the corrected one), which is why so many
# <v> shared-string indices shift in the raw XML diff even though only one
# cell's displayed text really changed.
#
# Non-content change: the sheet was scrolled/re-selected before saving
# (topLeftCell moved from A6 to A8, selection moved from B22 to G17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Fix the typo in the class description key ---
$ws.Range("B14").Value = "Grassland/Herbaceous"

# --- Reproduce the updated scroll position / active selection ---
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G17").Select()
